# Sentence caser: NLP adjustments and bug squashing
# Insert a new entry "Food and Drug Administration" into the alphabetised
# word list on Sheet1, column A. It belongs right before "GI Forum"
# (row 392), so insert a new cell there and push everything else down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 392

# Shift A392:A(end) down by one row, then write the new value into the
# freshly-opened cell.
$ws.Range("A" + $newRow).Insert(-4121)   # xlShiftDown
$ws.Range("A" + $newRow).Value = "Food and Drug Administration"

# Mirror the author's final selection/view state: the whole new row is
# selected (as if the row header for the inserted row had been clicked).
$ws.Rows($newRow).Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 378
    $win.ScrollColumn = 1
}
